$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 2 data rows (rows 2-3); the target has 4 data
# rows (rows 2-5). Insert two blank rows above row 2, pushing the old rows
# 2-3 down to rows 4-5. Clear formats on the freshly inserted rows so they
# don't inherit the bold/bordered header style that Insert() copies down
# from row 1.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).ClearFormats()

function Set-DataRow {
    param($ws, $r, $values)
    # Column B (Date) holds a string that looks like a date ("2025-10-10").
    # Assigning that directly to .Value lets Excel auto-convert it to a
    # date serial number; forcing a text NumberFormat first keeps it a
    # plain string, matching the source file.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
    # The text NumberFormat above is only needed to stop the date-like
    # string from being reinterpreted while it's being assigned. Reset the
    # cell style back to Normal afterwards so the saved file doesn't carry
    # a stray "@"-formatted style that the source file never had.
    $ws.Cells.Item($r, 2).Style = "Normal"
}

# --- Row 2: Croatian 2 HNL | Rudes vs BSK Bijelo Brdo ---
Set-DataRow $ws 2 @(
    "Croatian 2 HNL", "2025-10-10", "10:00:00", "Rudes", "BSK Bijelo Brdo", 1.02, 1000, 1.02, 1000, 1.02,
    950, 1.01, 1.05, 1.34, 1.05, 1.34, 1.32, 1.18, 1.32, 1.03,
    1.03, 1.01, 1.01, 1000, 1000, 1000, 1000, 1000, 1000, 1000,
    1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000,
    1000
)

# --- Row 3: Irish Premier Division | Shamrock Rovers vs Shelbourne ---
Set-DataRow $ws 3 @(
    "Irish Premier Division", "2025-10-10", "16:00:00", "Shamrock Rovers", "Shelbourne", 2.08, 2.26, 4, 4.7, 3.1,
    3.4, 0, 0, 0, 0, 1.54, 2.5, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0
)

# --- Row 4: FIFA World Cup Qualifiers - Americas | Bermuda vs Trinidad & Tobago ---
Set-DataRow $ws 4 @(
    "FIFA World Cup Qualifiers - Americas", "2025-10-10", "19:00:00", "Bermuda", "Trinidad & Tobago", 1.01, 1000, 1.01, 980, 1.01,
    1000, 0, 0, 0, 0, 1.24, 1.01, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0
)

# --- Row 5: FIFA World Cup Qualifiers - Americas | Curacao vs Jamaica ---
Set-DataRow $ws 5 @(
    "FIFA World Cup Qualifiers - Americas", "2025-10-10", "20:00:00", "Curacao", "Jamaica", 3.25, 4, 2.12, 2.46, 3.2,
    4.6, 0, 0, 0, 0, 1.7, 2.04, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0,
    0
)

$ws.Range("A1").Select()
